$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.190.68'
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").Value = '1.836.37'
$ws.Range("E3").Value = '  -0.16%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9995'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.26'
$ws.Range("E5").Value = '  +0.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6651'
$ws.Range("E6").Value = '  -3.03%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07373'
$ws.Range("E8").Value = '  -0.61%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2925'
$ws.Range("E9").Value = '  -2.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.68'
$ws.Range("E10").Value = '  -2.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07721'
$ws.Range("E11").Value = '  +1.05%  '
$ws.Range("D12").Value = '1.826.14'
$ws.Range("E12").Value = '  -0.77%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.980'
$ws.Range("E13").Value = '  -1.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6681'
$ws.Range("E14").Value = '  -1.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.76'
$ws.Range("E15").Value = '  -5.25%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.129'
$ws.Range("E16").Value = '  -0.40%  '
$ws.Range("D17").Value = '29.173.53'
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008275'
$ws.Range("E18").Value = '  +1.46%  '
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.44'
$ws.Range("E19").Value = '  -0.63%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '225.17'
$ws.Range("E20").Value = '  -1.82%  '
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.131'
$ws.Range("E22").Value = '  -3.07%  '
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '160.75'
$ws.Range("E24").Value = '  +0.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.627'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1393'
$ws.Range("E26").Value = '  -3.42%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.96'
$ws.Range("E27").Value = '  -0.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.508'
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.111'
$ws.Range("E29").Value = '  -3.55%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.032'
$ws.Range("E30").Value = '  -2.59%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.179'
$ws.Range("E31").Value = '  -1.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05303'
$ws.Range("E32").Value = '  +0.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.869'
$ws.Range("E33").Value = '  +0.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7531'
$ws.Range("E34").Value = '  -0.23%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.130'
$ws.Range("E35").Value = '  -0.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.676'
$ws.Range("E36").Value = '  -0.30%  '
$ws.Range("D37").Value = '1.296.55'
$ws.Range("E37").Value = '  +0.30%  '
$ws.Range("E38").Value = '  -1.83%  '
$ws.Range("E39").Value = '  +0.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9221'
$ws.Range("E40").Value = '  -1.93%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.08846'
$ws.Range("E41").Value = '  +18.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.956'
$ws.Range("E42").Value = '  +0.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.006'
$ws.Range("E43").Value = '  +0.66%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '102.19'
$ws.Range("E44").Value = '  -2.41%  '
$ws.Range("D45").Value = '1.976.86'
$ws.Range("E45").Value = '  -0.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5163'
$ws.Range("E46").Value = '  -0.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.767'
$ws.Range("E47").Value = '  +0.12%  '
$ws.Range("E48").Value = '  -1.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '63.07'
$ws.Range("E49").Value = '  -2.71%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05928'
$ws.Range("E50").Value = '  -0.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.019'
$ws.Range("E51").Value = '  -4.50%  '
